$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Silva / Astrid / dnc.silva.txt / Monday / speech
$ws.Range("A10").Value = "Silva"
$ws.Range("B10").Value = "Astrid"
$ws.Range("C10").Value = "dnc.silva.txt"
$ws.Range("D10").Value = "Monday"
$ws.Range("E10").Value = "speech"

# Row 16: Harkin / Tom / dnc.harkin.txt / Tuesday / speech
$ws.Range("A16").Value = "Harkin"
$ws.Range("B16").Value = "Tom"
$ws.Range("C16").Value = "dnc.harkin.txt"
$ws.Range("D16").Value = "Tuesday"
$ws.Range("E16").Value = "speech"

# Row 17: Grimes / Alison / dnc.grimes.txt / Tuesday / speech
$ws.Range("A17").Value = "Grimes"
$ws.Range("B17").Value = "Alison"
$ws.Range("C17").Value = "dnc.grimes.txt"
$ws.Range("D17").Value = "Tuesday"
$ws.Range("E17").Value = "speech"

# Update the view: scroll down and move the active selection
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A18").Select()
